$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "Classroom " -> "Name"
$ws.Range("A1").Value = "Name"

# Update selection to A2 (mirrors the saved view state in the target file)
$ws.Range("A2").Select()
